$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($ws, $rowNum, $values)
    $arr = New-Object "object[,]" 1,20
    for ($i = 0; $i -lt 20; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $rangeAddr = "A" + $rowNum + ":T" + $rowNum
    $ws.Range($rangeAddr).Value = $arr
}

Set-RowValues $ws 2 @("ECs", "Icam4", "Itgal", "ECs", 3, 1, 1.412252666666667, 4.236758, 0.409051039104018, 0.4090510391040179, 3, 1, 0.3547066666666667, 1.06412, 0.006934574868045491, 0.00693457486804549, 0.5009354358844444, 4.50841892296, 0.002836595055518617, 0.002836595055518616)
Set-RowValues $ws 3 @("ECs", "Icam4", "Itgal", "FAPs", 3, 1, 1.412252666666667, 4.236758, 0.409051039104018, 0.4090510391040179, 3, 1, 0.519977, 1.559931, 0.01016563762403213, 0.01016563762403213, 0.7343389048553334, 6.609050143698, 0.004158264633265244, 0.004158264633265242)
Set-RowValues $ws 4 @("ECs", "Icam4", "Itgal", "Inflammatory-Mac", 3, 1, 1.412252666666667, 4.236758, 0.409051039104018, 0.4090510391040179, 3, 1, 26.50170333333334, 79.50511, 0.5181127482682332, 0.5181127482682331, 37.42710120370889, 336.84391083338, 0.2119345580521593, 0.2119345580521592)
Set-RowValues $ws 5 @("ECs", "Icam4", "Itgal", "MuSCs", 3, 1, 1.412252666666667, 4.236758, 0.409051039104018, 0.4090510391040179, 3, 1, 0.09443866666666667, 0.283316, 0.001846291784117559, 0.001846291784117558, 0.1333712588364444, 1.200341329528, 0.0007552275727824985, 0.0007552275727824983)
Set-RowValues $ws 6 @("ECs", "Icam4", "Itgal", "Resolving-Mac", 3, 1, 1.412252666666667, 4.236758, 0.409051039104018, 0.4090510391040179, 3, 1, 23.67963033333333, 71.03889099999999, 0.4629407474555717, 0.4629407474555717, 33.44162108393088, 300.974589755378, 0.1893663937902924, 0.1893663937902923)
Set-RowValues $ws 7 @("FAPs", "Icam4", "Itgal", "ECs", 3, 1, 0.674352, 2.023056, 0.1953222626748136, 0.1953222626748136, 3, 1, 0.3547066666666667, 1.06412, 0.006934574868045491, 0.00693457486804549, 0.23919715008, 2.15277435072, 0.001354476853914543, 0.001354476853914542)
Set-RowValues $ws 8 @("FAPs", "Icam4", "Itgal", "FAPs", 3, 1, 0.674352, 2.023056, 0.1953222626748136, 0.1953222626748136, 3, 1, 0.519977, 1.559931, 0.01016563762403213, 0.01016563762403213, 0.350647529904, 3.155827769136, 0.001985575342258173, 0.001985575342258172)
Set-RowValues $ws 9 @("FAPs", "Icam4", "Itgal", "Inflammatory-Mac", 3, 1, 0.674352, 2.023056, 0.1953222626748136, 0.1953222626748136, 3, 1, 26.50170333333334, 79.50511, 0.5181127482682332, 0.5181127482682331, 17.87147664624, 160.84328981616, 0.1011989543124174, 0.1011989543124174)
Set-RowValues $ws 10 @("FAPs", "Icam4", "Itgal", "MuSCs", 3, 1, 0.674352, 2.023056, 0.1953222626748136, 0.1953222626748136, 3, 1, 0.09443866666666667, 0.283316, 0.001846291784117559, 0.001846291784117558, 0.063684903744, 0.5731641336960001, 0.0003606218888317601, 0.00036062188883176)
Set-RowValues $ws 11 @("FAPs", "Icam4", "Itgal", "Resolving-Mac", 3, 1, 0.674352, 2.023056, 0.1953222626748136, 0.1953222626748136, 3, 1, 23.67963033333333, 71.03889099999999, 0.4629407474555717, 0.4629407474555717, 15.968406074544, 143.715654670896, 0.09042263427739175, 0.09042263427739172)
Set-RowValues $ws 12 @("Inflammatory-Mac", "Icam4", "Itgal", "ECs", 3, 1, 0.9406819999999999, 2.822046, 0.2724632487150169, 0.2724632487150168, 3, 1, 0.3547066666666667, 1.06412, 0.006934574868045491, 0.00693457486804549, 0.3336661766133333, 3.00299558952, 0.001889416797005184, 0.001889416797005183)
Set-RowValues $ws 13 @("Inflammatory-Mac", "Icam4", "Itgal", "FAPs", 3, 1, 0.9406819999999999, 2.822046, 0.2724632487150169, 0.2724632487150168, 3, 1, 0.519977, 1.559931, 0.01016563762403213, 0.01016563762403213, 0.489133004314, 4.402197038825999, 0.0027697626523034, 0.002769762652303399)
Set-RowValues $ws 14 @("Inflammatory-Mac", "Icam4", "Itgal", "Inflammatory-Mac", 3, 1, 0.9406819999999999, 2.822046, 0.2724632487150169, 0.2724632487150168, 3, 1, 26.50170333333334, 79.50511, 0.5181127482682332, 0.5181127482682331, 24.92967529500667, 224.36707765506, 0.1411666825938285, 0.1411666825938285)
Set-RowValues $ws 15 @("Inflammatory-Mac", "Icam4", "Itgal", "MuSCs", 3, 1, 0.9406819999999999, 2.822046, 0.2724632487150169, 0.2724632487150168, 3, 1, 0.09443866666666667, 0.283316, 0.001846291784117559, 0.001846291784117558, 0.08883675383733333, 0.799530784536, 0.0005030466575765146, 0.0005030466575765144)
Set-RowValues $ws 16 @("Inflammatory-Mac", "Icam4", "Itgal", "Resolving-Mac", 3, 1, 0.9406819999999999, 2.822046, 0.2724632487150169, 0.2724632487150168, 3, 1, 23.67963033333333, 71.03889099999999, 0.4629407474555717, 0.4629407474555717, 22.27500202122066, 200.475018190986, 0.1261343400143033, 0.1261343400143032)
Set-RowValues $ws 17 @("MuSCs", "Icam4", "Itgal", "ECs", 1, 0.3333333333333333, 0.05847033333333334, 0.175411, 0.01693560307675702, 0.01693560307675701, 3, 1, 0.3547066666666667, 1.06412, 0.006934574868045491, 0.00693457486804549, 0.02073981703555556, 0.18665835332, 0.0001174412074712731, 0.0001174412074712731)
Set-RowValues $ws 18 @("MuSCs", "Icam4", "Itgal", "FAPs", 1, 0.3333333333333333, 0.05847033333333334, 0.175411, 0.01693560307675702, 0.01693560307675701, 3, 1, 0.519977, 1.559931, 0.01016563762403213, 0.01016563762403213, 0.03040322851566667, 0.273629056641, 0.0001721612038227555, 0.0001721612038227554)
Set-RowValues $ws 19 @("MuSCs", "Icam4", "Itgal", "Inflammatory-Mac", 1, 0.3333333333333333, 0.05847033333333334, 0.175411, 0.01693560307675702, 0.01693560307675701, 3, 1, 26.50170333333334, 79.50511, 0.5181127482682332, 0.5181127482682331, 1.549563427801111, 13.94607085021, 0.008774551853678524, 0.00877455185367852)
Set-RowValues $ws 20 @("MuSCs", "Icam4", "Itgal", "MuSCs", 1, 0.3333333333333333, 0.05847033333333334, 0.175411, 0.01693560307675702, 0.01693560307675701, 3, 1, 0.09443866666666667, 0.283316, 0.001846291784117559, 0.001846291784117558, 0.005521860319555556, 0.049696742876, [double]"3.126806481969253E-05", [double]"3.126806481969252E-05")
Set-RowValues $ws 21 @("MuSCs", "Icam4", "Itgal", "Resolving-Mac", 1, 0.3333333333333333, 0.05847033333333334, 0.175411, 0.01693560307675702, 0.01693560307675701, 3, 1, 23.67963033333333, 71.03889099999999, 0.4629407474555717, 0.4629407474555717, 1.384555878800111, 12.461002909201, 0.007840180746964773, 0.007840180746964772)
Set-RowValues $ws 22 @("Resolving-Mac", "Icam4", "Itgal", "ECs", 3, 1, 0.3667526666666667, 1.100258, 0.1062278464293945, 0.1062278464293945, 3, 1, 0.3547066666666667, 1.06412, 0.006934574868045491, 0.00693457486804549, 0.1300896158844445, 1.17080654296, 0.0007366449541358751, 0.0007366449541358749)
Set-RowValues $ws 23 @("Resolving-Mac", "Icam4", "Itgal", "FAPs", 3, 1, 0.3667526666666667, 1.100258, 0.1062278464293945, 0.1062278464293945, 3, 1, 0.519977, 1.559931, 0.01016563762403213, 0.01016563762403213, 0.1907029513553334, 1.716326562198, 0.00107987379238256, 0.00107987379238256)
Set-RowValues $ws 24 @("Resolving-Mac", "Icam4", "Itgal", "Inflammatory-Mac", 3, 1, 0.3667526666666667, 1.100258, 0.1062278464293945, 0.1062278464293945, 3, 1, 26.50170333333334, 79.50511, 0.5181127482682332, 0.5181127482682331, 9.71957036870889, 87.47613331837999, 0.05503800145614941, 0.05503800145614938)
Set-RowValues $ws 25 @("Resolving-Mac", "Icam4", "Itgal", "MuSCs", 3, 1, 0.3667526666666667, 1.100258, 0.1062278464293945, 0.1062278464293945, 3, 1, 0.09443866666666667, 0.283316, 0.001846291784117559, 0.001846291784117558, 0.03463563283644445, 0.311720695528, 0.0001961276001070928, 0.0001961276001070927)
Set-RowValues $ws 26 @("Resolving-Mac", "Icam4", "Itgal", "Resolving-Mac", 3, 1, 0.3667526666666667, 1.100258, 0.1062278464293945, 0.1062278464293945, 3, 1, 23.67963033333333, 71.03889099999999, 0.4629407474555717, 0.4629407474555717, 8.684567570430888, 78.16110813387799, 0.04917719862661958, 0.04917719862661956)
